# Generate Report for Handoff
# Renames the single tracked source file from "test-content-1.md" to "a.md",
# refreshes its handoff timestamps/target-file hashes, and adds a second
# tracked file "b.md" (duplicating the same handoff/target info) as a new
# row on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$HYPERLINK_COLOR = 15570276  # BGR int for RGB FF6495ED (matches existing "HyperLink" cell style)

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $HYPERLINK_COLOR
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 (existing "a.md", formerly "test-content-1.md")
$ov.Range("A2").Value = "a.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-31-13 14:31:06"
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
Style-AsHyperlink $ov.Range("A2")

# Row 3 (new "b.md")
$ov.Range("A3").Value = "b.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-31-13 14:31:06"
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null
Style-AsHyperlink $ov.Range("A3")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhHandoffDatetime = "2016-03-13 14:31:03"
$zhHandbackDatetime = "2016-03-13 08:34:50"

# Row 2
$zh.Range("A2").Value = "a.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = $zhXlf
$zh.Range("E2").Value = $zhHandoffDatetime
$zh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("F2").Value = "a.md"
$zh.Range("G2").Value = $zhXlf
$zh.Range("H2").Value = $zhHandbackDatetime
$zh.Range("I2").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/a.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46c22d5e14f08f4b581343d9f079c95343290fd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", [Type]::Missing, [Type]::Missing, $zhXlf) | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46c22d5e14f08f4b581343d9f079c95343290fd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", [Type]::Missing, [Type]::Missing, $zhXlf) | Out-Null

Style-AsHyperlink $zh.Range("A2")
Style-AsHyperlink $zh.Range("B2")
Style-AsHyperlink $zh.Range("D2")
Style-AsHyperlink $zh.Range("F2")
Style-AsHyperlink $zh.Range("G2")

# Row 3 (new "b.md")
$zh.Range("A3").Value = "b.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = $zhXlf
$zh.Range("E3").Value = $zhHandoffDatetime
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("F3").Value = "a.md"
$zh.Range("G3").Value = $zhXlf
$zh.Range("H3").Value = $zhHandbackDatetime
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/b.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46c22d5e14f08f4b581343d9f079c95343290fd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", [Type]::Missing, [Type]::Missing, $zhXlf) | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/74432397279876eb84635819af39f1bc8a0adcfd/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46c22d5e14f08f4b581343d9f079c95343290fd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", [Type]::Missing, [Type]::Missing, $zhXlf) | Out-Null

Style-AsHyperlink $zh.Range("A3")
Style-AsHyperlink $zh.Range("B3")
Style-AsHyperlink $zh.Range("D3")
Style-AsHyperlink $zh.Range("F3")
Style-AsHyperlink $zh.Range("G3")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deHandoffDatetime = "2016-03-13 14:31:06"
$deHandbackDatetime = "2016-03-13 08:34:56"

# Row 2
$de.Range("A2").Value = "a.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = $deXlf
$de.Range("E2").Value = $deHandoffDatetime
$de.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("F2").Value = "a.md"
$de.Range("G2").Value = $deXlf
$de.Range("H2").Value = $deHandbackDatetime
$de.Range("I2").Value = "Include"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c9a887eca8d06e1f356440d67f26d32c46de717/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c9a887eca8d06e1f356440d67f26d32c46de717/e2e/a.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2ced4a799e6cc7df05ec2566674f65fda088143/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", [Type]::Missing, [Type]::Missing, $deXlf) | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c9a887eca8d06e1f356440d67f26d32c46de717/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2ced4a799e6cc7df05ec2566674f65fda088143/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", [Type]::Missing, [Type]::Missing, $deXlf) | Out-Null

Style-AsHyperlink $de.Range("A2")
Style-AsHyperlink $de.Range("B2")
Style-AsHyperlink $de.Range("D2")
Style-AsHyperlink $de.Range("F2")
Style-AsHyperlink $de.Range("G2")

# Row 3 (new "b.md")
$de.Range("A3").Value = "b.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = $deXlf
$de.Range("E3").Value = $deHandoffDatetime
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("F3").Value = "a.md"
$de.Range("G3").Value = $deXlf
$de.Range("H3").Value = $deHandbackDatetime
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c9a887eca8d06e1f356440d67f26d32c46de717/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c9a887eca8d06e1f356440d67f26d32c46de717/e2e/b.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2ced4a799e6cc7df05ec2566674f65fda088143/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", [Type]::Missing, [Type]::Missing, $deXlf) | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c9a887eca8d06e1f356440d67f26d32c46de717/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2ced4a799e6cc7df05ec2566674f65fda088143/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", [Type]::Missing, [Type]::Missing, $deXlf) | Out-Null

Style-AsHyperlink $de.Range("A3")
Style-AsHyperlink $de.Range("B3")
Style-AsHyperlink $de.Range("D3")
Style-AsHyperlink $de.Range("F3")
Style-AsHyperlink $de.Range("G3")

"Report regenerated: a.md refreshed, b.md added to Overview/zh-cn/de-de."
